$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update capital structure database values for rows 2 and 3 (identical values in both rows)
$ws.Range("D2:D3").Value = 0.0445
$ws.Range("E2:E3").Value = -0.0169
$ws.Range("G2:G3").Value = 0.5163398692810458
$ws.Range("H2:H3").Value = 0.5163398692810458
$ws.Range("I2:I3").Value = 0.4640522875816994
$ws.Range("J2:J3").Value = 0.4640522875816994
$ws.Range("K2:K3").Value = 19.6
$ws.Range("L2:L3").Value = 0.4270152505446623
$ws.Range("M2:M3").Value = 19.8
$ws.Range("N2:N3").Value = 0.06964474147027787
$ws.Range("O2:O3").Value = 1.010204081632653
$ws.Range("P2:P3").Value = 19.8
$ws.Range("Q2:Q3").Value = 0.06964474147027787
$ws.Range("R2:R3").Value = 1.010204081632653
$ws.Range("U2:U3").Value = 80.09999999999999
$ws.Range("V2:V3").Value = 0.2817446359479423
$ws.Range("W2:W3").Value = 0.1946375372393247
$ws.Range("X2:X3").Value = 0.04436545039956676
$ws.Range("Y2:Y3").Value = 0.150272086839758
$ws.Range("Z2:Z3").Value = 1.599303135888501
$ws.Range("AA2:AA3").Value = 0.7421602787456445
$ws.Range("AB2:AB3").Value = 0.04436545039956676
$ws.Range("AC2:AC3").Value = 0.6977948283460778
$ws.Range("AG2:AG3").Value = -80.09999999999999
$ws.Range("AJ2:AJ3").Value = -0.3922624877571008
$ws.Range("AK2:AK3").Value = -3.282786885245901
$ws.Range("AP2:AP3").Value = -3.467532467532467
